$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert four new bullet paragraphs right after the "JavaScript" heading
#    and before "Collapse function for FAQS and Get Involved page":
#      - Sticky function for navigation bar            (ilvl 0)
#      - If user goes past the position of the navigation bar  (ilvl 1)
#      - Call the Sticky function                       (ilvl 2)
#      - Else don't call the sticky function             (ilvl 1)
# ---------------------------------------------------------------------------

$anchor = $d.Content.Find
$found = $d.Content.Find.Execute("Collapse function for FAQS and Get Involved page",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Locate the paragraph that holds "Collapse function for FAQS and Get Involved page"
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Collapse function for FAQS and Get Involved page`r") {
        $targetIndex = $i
        break
    }
}

# Paragraph 1: "Sticky function for navigation bar " (ilvl 0, same indent as target)
$d.Paragraphs.Item($targetIndex).Range.InsertParagraphBefore()
$targetIndex = $targetIndex + 1
$p1 = $d.Paragraphs.Item($targetIndex - 1)
$p1.Range.Text = "Sticky function for navigation bar "

# Paragraph 2: "If user goes past the position of the navigation bar" (ilvl 1)
$d.Paragraphs.Item($targetIndex).Range.InsertParagraphBefore()
$targetIndex = $targetIndex + 1
$p2 = $d.Paragraphs.Item($targetIndex - 1)
$p2.Range.ListFormat.ListLevelNumber = 2
$p2.Format.LeftIndent = 72
$p2.Format.FirstLineIndent = -18
$p2.Range.Text = "If user goes past the position of the navigation bar"

# Paragraph 3: "Call the Sticky function" (ilvl 2)
$d.Paragraphs.Item($targetIndex).Range.InsertParagraphBefore()
$targetIndex = $targetIndex + 1
$p3 = $d.Paragraphs.Item($targetIndex - 1)
$p3.Range.ListFormat.ListLevelNumber = 3
$p3.Format.LeftIndent = 108
$p3.Format.FirstLineIndent = -18
$p3.Range.Text = "Call the Sticky function"

# Paragraph 4: "Else don't call the sticky function" (ilvl 1)
$d.Paragraphs.Item($targetIndex).Range.InsertParagraphBefore()
$targetIndex = $targetIndex + 1
$p4 = $d.Paragraphs.Item($targetIndex - 1)
$p4.Range.ListFormat.ListLevelNumber = 2
$p4.Format.LeftIndent = 72
$p4.Format.FirstLineIndent = -18
$p4.Range.Text = "Else don't call the sticky function"

# ---------------------------------------------------------------------------
# 2. Remove the old "Drop Down menu function for navigation bar " paragraph
#    and the blank bullet paragraph immediately following it (they used to
#    sit right after "Sets the maxHeight to be null or to have some pixels").
# ---------------------------------------------------------------------------

$dropIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Drop Down menu function for navigation bar `r") {
        $dropIndex = $i
        break
    }
}

# Delete the "Drop Down menu function for navigation bar " paragraph entirely
$d.Paragraphs.Item($dropIndex).Range.Delete()
# The following blank bullet paragraph now sits at the same index; delete it too
$d.Paragraphs.Item($dropIndex).Range.Delete()
